$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural change: two new BOM lines were added (one per assembly) ---
# Insert a new row at position 10 (front assembly gains a "Spacer" line),
# which shifts the old rows 10-18 down to 11-19.
$ws.Rows.Item(10).Insert()

# Insert a new row at position 20 (rear assembly gains an "Aluminium tapped tube" line).
$ws.Rows.Item(20).Insert()

# Give the two freshly inserted rows the same formatting (borders/fill/font) as the
# other data rows in their respective blocks, instead of Excel's default blank style.
$ws.Range("A9:G9").Copy()
$ws.Range("A10:G10").PasteSpecial(-4122)
$ws.Range("A19:G19").Copy()
$ws.Range("A20:G20").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Front assembly (SU_A0014) content updates ---

# Row 4: Right Bearing Support - now documents which part its diameter matches
$ws.Cells.Item(4,5).Value = "Outside diameter of the SU_14007"

# Row 5: Left Bearing Support - same comment added
$ws.Cells.Item(5,5).Value = "Outside diameter of the SU_14007"

# Row 7: the combined "Rod ends bearing" (qty 4) line is split into a male-only line (qty 2)
$ws.Cells.Item(7,3).Value = "Rod ends bearing, male"
$ws.Cells.Item(7,4).Value = "b"
$ws.Cells.Item(7,5).Value = "2 male thread, ARB rod, right hand thread"
$ws.Cells.Item(7,6).Value = 2
$ws.Cells.Item(7,7).Value = "SU_14005"

# Row 8: new "Rod ends bearing, female" line (qty 2)
$ws.Cells.Item(8,3).Value = "Rod ends bearing, female"
$ws.Cells.Item(8,4).Value = "b"
$ws.Cells.Item(8,5).Value = "2 female thread, ARB rod, right hand thread"
$ws.Cells.Item(8,6).Value = 2
$ws.Cells.Item(8,7).Value = "SU_14006"

# Row 9: "Spherical plain bearings" line (shifted down one), ID renumbered to SU_14007
$ws.Cells.Item(9,3).Value = "Spherical plain bearings"
$ws.Cells.Item(9,4).Value = "b"
$ws.Cells.Item(9,5).Value = "Used in the bearing supports"
$ws.Cells.Item(9,6).Value = 2
$ws.Cells.Item(9,7).Value = "SU_14007"

# Row 10 (new): "Spacer" line, ID SU_14008
$ws.Cells.Item(10,3).Value = "Spacer "
$ws.Cells.Item(10,4).Value = "m"
$ws.Cells.Item(10,5).Value = "M6 type 16 mm spacer"
$ws.Cells.Item(10,6).Value = 4
$ws.Cells.Item(10,7).Value = "SU_14008"

# --- Rear assembly (SU_A0015) content updates ---

# Row 13: Right Bearing Support - comment added
$ws.Cells.Item(13,5).Value = "Outside diameter of the SU_15009"

# Row 14: Left Bearing Support - comment added
$ws.Cells.Item(14,5).Value = "Outside diameter of the SU_15009"

# Row 17: the combined "Rod ends bearing" (qty 4) line becomes the left-hand-thread male line (qty 2)
$ws.Cells.Item(17,3).Value = "Rod ends bearing, male"
$ws.Cells.Item(17,4).Value = "b"
$ws.Cells.Item(17,5).Value = "2 with a left-hand thread, ARB rod"
$ws.Cells.Item(17,6).Value = 2
$ws.Cells.Item(17,7).Value = "SU_15006"

# Row 18: new right-hand-thread male line
$ws.Cells.Item(18,3).Value = "Rod ends bearing, male"
$ws.Cells.Item(18,4).Value = "b"
$ws.Cells.Item(18,5).Value = "2 with a right-hand thread, ARB rod "
$ws.Cells.Item(18,6).Value = 2
$ws.Cells.Item(18,7).Value = "SU_15006"

# Row 19: "Spherical plain bearings" line (shifted down), ID SU_15008
$ws.Cells.Item(19,3).Value = "Spherical plain bearings"
$ws.Cells.Item(19,4).Value = "b"
$ws.Cells.Item(19,5).Value = "Used in the bearing supports"
$ws.Cells.Item(19,6).Value = 2
$ws.Cells.Item(19,7).Value = "SU_15008"

# Row 20 (new): "Aluminium tapped tube" line (replaces plain "Aluminium tube" wording), ID SU_15009
$ws.Cells.Item(20,3).Value = "Aluminium tapped tube"
$ws.Cells.Item(20,4).Value = "m"
$ws.Cells.Item(20,5).Value = "Aluminium tube for the ARB rod"
$ws.Cells.Item(20,6).Value = 2
$ws.Cells.Item(20,7).Value = "SU_15009"
